$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds free-form price text that can look numeric (e.g. "213.29"),
# and must stay plain text exactly like the source OOXML (inlineStr, no thousands
# grouping, no style change). Force text entry via NumberFormat "@", then restore
# the original (unstyled) look with the Normal cell style so no new formatting
# gets attached to the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "28.428.40"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.583.73"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "213.29"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.492"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
Set-TextValue $ws.Range("D8") "44.53"
$ws.Range("E8").Value = "  -1.99%  "

# Row 9
Set-TextValue $ws.Range("D9") "23.92"
$ws.Range("E9").Value = "  -1.22%  "

# Row 10
$ws.Range("E10").Value = "  -1.69%  "

# Row 11
$ws.Range("E11").Value = "  -1.72%  "

# Row 12
$ws.Range("E12").Value = "  +1.06%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.810.09"
$ws.Range("E13").Value = "  -0.19%  "

# Row 14
Set-TextValue $ws.Range("D14") "1.581.57"
$ws.Range("E14").Value = "  -0.78%  "

# Row 15
$ws.Range("E15").Value = "  -0.98%  "

# Row 16
$ws.Range("E16").Value = "  -1.85%  "

# Row 17
Set-TextValue $ws.Range("D17") "28.437.98"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
Set-TextValue $ws.Range("D18") "62.04"
$ws.Range("E18").Value = "  -1.65%  "

# Row 19
Set-TextValue $ws.Range("D19") "229.95"
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.47"
$ws.Range("E20").Value = "  -0.41%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.0₃0688"
$ws.Range("E21").Value = "  -2.30%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("E23").Value = "  -3.27%  "

# Row 24
Set-TextValue $ws.Range("D24") "9.16"
$ws.Range("E24").Value = "  -1.77%  "

# Row 25
$ws.Range("E25").Value = "  +3.19%  "

# Row 26
Set-TextValue $ws.Range("D26") "151.65"
$ws.Range("E26").Value = "  -0.17%  "

# Row 27
Set-TextValue $ws.Range("D27") "15.05"
$ws.Range("E27").Value = "  -1.14%  "

# Row 28
$ws.Range("E28").Value = "  -1.79%  "

# Row 29
$ws.Range("E29").Value = "  -1.75%  "

# Row 30
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("E31").Value = "  +2.45%  "

# Row 32
$ws.Range("E32").Value = "  -1.24%  "

# Row 33
$ws.Range("E33").Value = "  -1.48%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.09"
$ws.Range("E34").Value = "  -2.43%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.395.66"
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("E36").Value = "  +7.29%  "

# Row 37
$ws.Range("E37").Value = "  -5.13%  "

# Row 38
$ws.Range("E38").Value = "  +0.19%  "

# Row 39
$ws.Range("E39").Value = "  +1.14%  "

# Row 40
$ws.Range("E40").Value = "  -1.05%  "

# Row 41
$ws.Range("E41").Value = "  -3.28%  "

# Row 42
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("E43").Value = "  -2.55%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.89"
$ws.Range("E44").Value = "  +1.00%  "

# Row 45
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D45") "0.0460"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D46") "5.44"
$ws.Range("E46").Value = "  -3.53%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.944"
$ws.Range("E47").Value = "  -3.88%  "

# Row 48
Set-TextValue $ws.Range("D48") "62.81"
$ws.Range("E48").Value = "  -0.18%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.721.28"
$ws.Range("E49").Value = "  -0.12%  "

# Row 50
Set-TextValue $ws.Range("D50") "86.69"
$ws.Range("E50").Value = "  +0.09%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0₆0103"
$ws.Range("E51").Value = "  -0.99%  "
